$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7880238890647888
$ws.Range("B1").Value = 1.480337977409363
$ws.Range("C1").Value = 5.703042030334473
$ws.Range("D1").Value = 3.147686004638672
$ws.Range("E1").Value = 1.48246967792511
